$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" "69.502.48"
$ws.Range("E2").Value = "  -0.17%  "
Set-TextValue "D3" "3.780.25"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "615.69"
$ws.Range("E5").Value = "  -0.09%  "
Set-TextValue "D6" "177.46"
$ws.Range("E6").Value = "  +0.40%  "
Set-TextValue "D7" "3.778.97"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("E8").Value = "  +0.08%  "
Set-TextValue "D9" "0.528"
$ws.Range("E9").Value = "  -0.63%  "
Set-TextValue "D10" "0.166"
$ws.Range("E10").Value = "  -0.37%  "
Set-TextValue "D11" "6.51"
$ws.Range("E11").Value = "  +3.99%  "
Set-TextValue "D12" "0.486"
$ws.Range("E12").Value = "  +0.17%  "
Set-TextValue "D13" "39.87"
$ws.Range("E13").Value = "  -2.18%  "
Set-TextValue "D14" "0.0000254"
$ws.Range("E14").Value = "  -0.31%  "
Set-TextValue "D15" "4.415.89"
$ws.Range("E15").Value = "  +1.31%  "
Set-TextValue "D16" "3.783.92"
$ws.Range("E16").Value = "  +1.06%  "
Set-TextValue "D17" "69.583.85"
$ws.Range("E17").Value = "  -0.10%  "
Set-TextValue "D18" "7.57"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  -3.16%  "
Set-TextValue "D20" "509.24"
$ws.Range("E20").Value = "  +0.30%  "
Set-TextValue "D21" "16.47"
$ws.Range("E21").Value = "  -0.33%  "
Set-TextValue "D22" "9.42"
Set-TextValue "D23" "0.734"
$ws.Range("E23").Value = "  +2.01%  "
Set-TextValue "D24" "2.49"
$ws.Range("E24").Value = "  +0.04%  "
Set-TextValue "D25" "86.25"
$ws.Range("E25").Value = "  -0.32%  "
Set-TextValue "D26" "12.89"
$ws.Range("E26").Value = "  -1.43%  "
Set-TextValue "D27" "0.0000139"
$ws.Range("E27").Value = "  +3.23%  "
Set-TextValue "D28" "10.56"
$ws.Range("E28").Value = "  -3.93%  "
$ws.Range("E29").Value = "  +0.46%  "
Set-TextValue "D30" "2.53"
$ws.Range("E30").Value = "  +2.27%  "
Set-TextValue "D31" "3.00"
$ws.Range("E31").Value = "  +4.20%  "
Set-TextValue "D32" "8.07"
$ws.Range("E32").Value = "  +3.71%  "
Set-TextValue "D33" "31.01"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  -0.02%  "
Set-TextValue "D36" "1.05"
$ws.Range("E36").Value = "  -1.21%  "
Set-TextValue "D37" "6.13"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  +7.82%  "
Set-TextValue "D39" "0.340"
$ws.Range("E39").Value = "  +2.06%  "
Set-TextValue "D40" "462.22"
$ws.Range("E40").Value = "  +10.01%  "
$ws.Range("E41").Value = "  -1.51%  "
Set-TextValue "D42" "49.80"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("E43").Value = "  +8.81%  "
Set-TextValue "D44" "44.34"
$ws.Range("E44").Value = "  -0.59%  "
Set-TextValue "D45" "8.57"
$ws.Range("E45").Value = "  -0.74%  "
Set-TextValue "D46" "2.960.77"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("E47").Value = "  +0.96%  "
Set-TextValue "D48" "27.36"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +0.04%  "
Set-TextValue "D50" "138.54"
$ws.Range("E50").Value = "  -0.42%  "
Set-TextValue "D51" "2.47"
$ws.Range("E51").Value = "  -0.64%  "
